# Rename the embedded logo pictures' shape names.
#
#   Footer(1)  (primary footer)    Pearson logo  id=1 : image2.png -> image1.png
#   Footer(2)  (first-page footer) Pearson logo  id=2 : image2.png -> image1.png
#   Header(2)  (first-page header) BTEC logo     id=3 : image1.jpg -> image2.jpg
#
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument

function Rename-LogoShape($range, [string]$oldName, [string]$newName) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -eq $oldName) {
            $shape.Name = $newName
        }
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $section = $d.Sections.Item($si)

    $footer1 = $section.Footers.Item(1)
    if ($footer1.Exists) {
        Rename-LogoShape $footer1.Range "image2.png" "image1.png"
    }

    $footer2 = $section.Footers.Item(2)
    if ($footer2.Exists) {
        Rename-LogoShape $footer2.Range "image2.png" "image1.png"
    }

    $header2 = $section.Headers.Item(2)
    if ($header2.Exists) {
        Rename-LogoShape $header2.Range "image1.jpg" "image2.jpg"
    }
}
